# Tukey's post-hoc of CAZyme domain data
# Adds a new "CAZyme domains" worksheet (after litterChemistry) containing
# the Tukey post-hoc significance table for the CAZyme domain substrates,
# and updates the previously-active sheet's selection/tab state.

$wb = $excel.ActiveWorkbook

# --- Deselect / move off the old active sheet (litterChemistry) first, so
#     that the new sheet ends up as the one and only tabSelected sheet. ---
$ws3 = $wb.Worksheets.Item("litterChemistry")
$ws3.Range("C10").Select() | Out-Null

# --- Add the new worksheet as the last tab in the workbook. ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CAZyme domains"

# --- Column widths (best-fit in the source file). ---
$ws.Columns.Item(1).ColumnWidth = 13.94
$ws.Columns.Item(2).ColumnWidth = 8.39
$ws.Columns.Item(3).ColumnWidth = 9.39
$ws.Columns.Item(4).ColumnWidth = 10.83
$ws.Columns.Item(5).ColumnWidth = 21.17
$ws.Columns.Item(6).ColumnWidth = 19.61
$ws.Columns.Item(7).ColumnWidth = 22.17
$ws.Columns.Item(8).ColumnWidth = 9.05

# --- Table data: header row + 14 substrate rows across 8 columns. ---
$data = New-Object 'object[,]' 15,8

$data[0,0] = "Substrate"
$data[0,1] = "timePoint"
$data[0,2] = "Vegetation"
$data[0,3] = "Precipitation"
$data[0,4] = "timePoint x Precipitation"
$data[0,5] = "timePoint x Vegetation"
$data[0,6] = "Vegetation x Precipitation"
$data[0,7] = "Three-way"

$data[1,0] = "Hemicellulose";   $data[1,1] = "***"; $data[1,2] = "***"; $data[1,3] = "o";   $data[1,4] = "o"; $data[1,5] = "*";   $data[1,6] = "o";   $data[1,7] = "o"
$data[2,0] = "Lignin";          $data[2,1] = "*";   $data[2,2] = "***"; $data[2,3] = "***"; $data[2,4] = "**"; $data[2,5] = "o";   $data[2,6] = "o";   $data[2,7] = "*"
$data[3,0] = "Polysaccharide";  $data[3,1] = "***"; $data[3,2] = "***"; $data[3,3] = "o";   $data[3,4] = "o"; $data[3,5] = "o";   $data[3,6] = "o";   $data[3,7] = "*"
$data[4,0] = "Oligosaccharides";$data[4,1] = "*";   $data[4,2] = "***"; $data[4,3] = "o";   $data[4,4] = "o"; $data[4,5] = "***"; $data[4,6] = "o";   $data[4,7] = "o"
$data[5,0] = "Cell_wall";       $data[5,1] = "***"; $data[5,2] = "***"; $data[5,3] = "o";   $data[5,4] = "o"; $data[5,5] = "o";   $data[5,6] = "o";   $data[5,7] = "*"
$data[6,0] = "Inulin";          $data[6,1] = "***"; $data[6,2] = "*";   $data[6,3] = "o";   $data[6,4] = "o"; $data[6,5] = "o";   $data[6,6] = "***"; $data[6,7] = "**"
$data[7,0] = "Starch";          $data[7,1] = "o";   $data[7,2] = "o";   $data[7,3] = "o";   $data[7,4] = "o"; $data[7,5] = "o";   $data[7,6] = "*";   $data[7,7] = "o"
$data[8,0] = "Trehalose";       $data[8,1] = "***"; $data[8,2] = "**";  $data[8,3] = "o";   $data[8,4] = "o"; $data[8,5] = "o";   $data[8,6] = "o";   $data[8,7] = "o"
$data[9,0] = "Cellulose";       $data[9,1] = "***"; $data[9,2] = "**";  $data[9,3] = "o";   $data[9,4] = "o"; $data[9,5] = "o";   $data[9,6] = "***"; $data[9,7] = "***"
$data[10,0] = "Pectin";         $data[10,1] = "**"; $data[10,2] = "o";  $data[10,3] = "**"; $data[10,4] = "o"; $data[10,5] = "o"; $data[10,6] = "***"; $data[10,7] = "o"
$data[11,0] = "Glycogen";       $data[11,1] = "**"; $data[11,2] = "o";  $data[11,3] = "o";  $data[11,4] = "o"; $data[11,5] = "**"; $data[11,6] = "o";   $data[11,7] = "o"
$data[12,0] = "Peptidoglycan";  $data[12,1] = "**"; $data[12,2] = "o";  $data[12,3] = "o";  $data[12,4] = "o"; $data[12,5] = "o";  $data[12,6] = "o";   $data[12,7] = "o"
$data[13,0] = "Chitin";         $data[13,1] = "**"; $data[13,2] = "o";  $data[13,3] = "o";  $data[13,4] = "o"; $data[13,5] = "o";  $data[13,6] = "o";   $data[13,7] = "o"
$data[14,0] = "Total";          $data[14,1] = "**"; $data[14,2] = "***";$data[14,3] = "o";  $data[14,4] = "o"; $data[14,5] = "*";  $data[14,6] = "o";   $data[14,7] = "*"

$ws.Range("A1:H15").Value = $data

# --- Formatting: bold header row with borders, and borders on the data rows. ---
$headerRange = $ws.Range("A1:H1")
$headerRange.Font.Bold = $true
$headerRange.Borders.ColorIndex = 1
$headerRange.Borders.LineStyle = 1

$bodyRange = $ws.Range("A2:H15")
$bodyRange.Borders.ColorIndex = 1
$bodyRange.Borders.LineStyle = 1

# --- Final selection on the new sheet, which also makes it the active tab. ---
$ws.Range("B18").Select() | Out-Null
